# edit.ps1 - applies the four text edits described by the diff:
#  1. Slide 27: merge the "(5 x 24) + (24 x 2) + 24 + 2 = 194!" runs into one run.
#  2. Slide 29: fix "Fell" -> "Feel" and split off a leading space run.
#  3. Slide 5 : merge two pairs of split runs back into single runs (no text change).
#  4. Slide 8 : merge " " + "input" runs into a single " input" run.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 27 - "(5 x 24) + (24 x 2) + 24 + 2 = 194!"
# ---------------------------------------------------------------------------
$sh27 = $p.Slides.Item(27).Shapes.Item(2)
$tr27 = $sh27.TextFrame.TextRange
$full27 = $tr27.Text
$newText27 = "(5 x 24) + (24 x 2) + 24 + 2 = 194!"
$idx27 = $full27.IndexOf("(5 x 24)")
$c27 = $tr27.Characters($idx27 + 1, $newText27.Length)
$c27.Text = $newText27

# ---------------------------------------------------------------------------
# 2) Slide 29 - "Fell" -> "Feel" / " " / "free to test and "
# ---------------------------------------------------------------------------
$sh29 = $p.Slides.Item(29).Shapes.Item(2)
$tr29 = $sh29.TextFrame.TextRange
$full29 = $tr29.Text
$idx29 = $full29.IndexOf("Fell")

# "Fell" -> "Feel"
$cWord = $tr29.Characters($idx29 + 1, 4)
$cWord.Text = "Feel"

# the character right after "Feel" is the space - rewrite it on its own so it
# becomes its own run
$tr29b = $sh29.TextFrame.TextRange
$cSpace = $tr29b.Characters($idx29 + 4 + 1, 1)
$cSpace.Text = " "

# the remaining "free to test and " text, rewritten as its own run
$rest = "free to test and "
$tr29c = $sh29.TextFrame.TextRange
$cRest = $tr29c.Characters($idx29 + 4 + 1 + 1, $rest.Length)
$cRest.Text = $rest

# ---------------------------------------------------------------------------
# 3) Slide 5 - merge split runs (text itself is unchanged)
# ---------------------------------------------------------------------------
$sh5 = $p.Slides.Item(5).Shapes.Item(2)
$tr5 = $sh5.TextFrame.TextRange
$full5 = $tr5.Text

$seg1 = 'For our purposes, we can simplified this action in a way such that an "importance factor" called '
$idx5a = $full5.IndexOf("For our purposes")
$c5a = $tr5.Characters($idx5a + 1, $seg1.Length)
$c5a.Text = $seg1

$tr5b = $sh5.TextFrame.TextRange
$full5b = $tr5b.Text
$seg2 = ' is assigned to each input. The neuron will transmit a single value that belongs to '
$idx5b = $full5b.IndexOf(" is assigned to each input. The neuron will tra")
$c5b = $tr5b.Characters($idx5b + 1, $seg2.Length)
$c5b.Text = $seg2

# ---------------------------------------------------------------------------
# 4) Slide 8 - merge " " + "input" into a single " input" run
# ---------------------------------------------------------------------------
$sh8 = $p.Slides.Item(8).Shapes.Item(2)
$tr8 = $sh8.TextFrame.TextRange
$full8 = $tr8.Text
$anchor8 = $full8.IndexOf("wi is the weight of the ith input")
$prefix8 = "wi is the weight of the ith"
$target8 = $anchor8 + $prefix8.Length
$c8 = $tr8.Characters($target8 + 1, 6)
$c8.Text = " input"
